$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.674.19'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.646.43'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.23%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.06'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  -0.50%  '
$ws.Range('E9').Value = '  +0.37%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.37'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  +0.35%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0843'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').Value = '1.877.24'
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('E13').Value = '  +2.97%  '
$ws.Range('D14').Value = '1.643.15'
$ws.Range('E14').Value = '  +0.33%  '
$style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.535'
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +1.42%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.33'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +4.26%  '
$ws.Range('D17').Value = '26.736.90'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '0.0₃0755'
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  +0.24%  '
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.39'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  +1.55%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.33'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +1.94%  '
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.56'
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('E24').Value = '  +10.56%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.30'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  +2.59%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.93'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +2.45%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0518'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('E32').Value = '  +2.24%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.07'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +2.49%  '
$ws.Range('D34').Value = '1.286.24'
$ws.Range('E34').Value = '  +5.44%  '
$ws.Range('E35').Value = '  +1.26%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0184'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  +6.18%  '
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('E38').Value = '  +3.93%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.828'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +1.70%  '
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('E42').Value = '  -2.11%  '
$style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.44'
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').Value = '1.788.98'
$ws.Range('E44').Value = '  +0.83%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.81'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  +0.15%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.14'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +9.33%  '
$ws.Range('E47').Value = '  +3.65%  '
$ws.Range('E48').Value = '  +0.51%  '
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.84'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +1.45%  '
$style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0979'
$ws.Range('D50').Style = $style
$ws.Range('E50').Value = '  +3.07%  '
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.408'
$ws.Range('D51').Style = $style
$ws.Range('E51').Value = '  -0.71%  '
